$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay as text (matches source formatting)
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply updated price (D) and volume-change (E) values
$ws.Range("D2").Value = "28.372.30"
$ws.Range("E2").Value = "  +0.35%  "
$ws.Range("D3").Value = "1.867.39"
$ws.Range("E3").Value = "  -0.15%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "330.40"
$ws.Range("E5").Value = "  -2.16%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("D7").Value = "0.4588"
$ws.Range("E7").Value = "  -2.39%  "
$ws.Range("D8").Value = "0.4001"
$ws.Range("E8").Value = "  +1.80%  "
$ws.Range("D9").Value = "47.57"
$ws.Range("E9").Value = "  +0.85%  "
$ws.Range("D10").Value = "0.07838"
$ws.Range("E10").Value = "  -1.85%  "
$ws.Range("D11").Value = "0.9832"
$ws.Range("E11").Value = "  -1.93%  "
$ws.Range("D12").Value = "21.24"
$ws.Range("E12").Value = "  -2.26%  "
$ws.Range("D13").Value = "1.881.63"
$ws.Range("E13").Value = "  -0.20%  "
$ws.Range("D14").Value = "5.837"
$ws.Range("E14").Value = "  -2.52%  "
$ws.Range("D15").Value = "6.983"
$ws.Range("E15").Value = "  -4.05%  "
$ws.Range("D16").Value = "1.002"
$ws.Range("E16").Value = "  -0.01%  "
$ws.Range("D17").Value = "88.08"
$ws.Range("E17").Value = "  -3.37%  "
$ws.Range("D18").Value = "0.06532"
$ws.Range("E18").Value = "  -0.76%  "
$ws.Range("D19").Value = "0.00001016"
$ws.Range("E19").Value = "  -2.64%  "
$ws.Range("D20").Value = "17.12"
$ws.Range("E20").Value = "  -3.12%  "
$ws.Range("E21").Value = "  +0.23%  "
$ws.Range("D22").Value = "28.350.33"
$ws.Range("E22").Value = "  +0.22%  "
$ws.Range("D23").Value = "5.325"
$ws.Range("D24").Value = "10.83"
$ws.Range("E24").Value = "  -1.77%  "
$ws.Range("E25").Value = "  -2.00%  "
$ws.Range("D26").Value = "2.099.71"
$ws.Range("E26").Value = "  -0.38%  "
$ws.Range("D27").Value = "157.47"
$ws.Range("E27").Value = "  -1.14%  "
$ws.Range("D28").Value = "19.30"
$ws.Range("E28").Value = "  -2.76%  "
$ws.Range("D29").Value = "2.055"
$ws.Range("E29").Value = "  -4.84%  "
$ws.Range("D30").Value = "5.280"
$ws.Range("E30").Value = "  -3.98%  "
$ws.Range("D31").Value = "117.12"
$ws.Range("E31").Value = "  -2.34%  "
$ws.Range("D32").Value = "0.9521"
$ws.Range("E32").Value = "  -2.80%  "
$ws.Range("D33").Value = "0.09309"
$ws.Range("E33").Value = "  -1.85%  "
$ws.Range("D34").Value = "3.593"
$ws.Range("E34").Value = "  +0.30%  "
$ws.Range("D35").Value = "1.382"
$ws.Range("E35").Value = "  +0.11%  "
$ws.Range("D36").Value = "5.220"
$ws.Range("E36").Value = "  -2.56%  "
$ws.Range("D37").Value = "0.06013"
$ws.Range("E37").Value = "  -1.16%  "
$ws.Range("D38").Value = "0.02198"
$ws.Range("E38").Value = "  -3.33%  "
$ws.Range("D39").Value = "8.268"
$ws.Range("E39").Value = "  -2.27%  "
$ws.Range("D40").Value = "1.163"
$ws.Range("E40").Value = "  -1.51%  "
$ws.Range("D41").Value = "1.001"
$ws.Range("E41").Value = "  +0.08%  "
$ws.Range("D42").Value = "0.5742"
$ws.Range("E42").Value = "  -3.70%  "
$ws.Range("D43").Value = "0.1804"
$ws.Range("E43").Value = "  -3.86%  "
$ws.Range("D44").Value = "9.988"
$ws.Range("E44").Value = "  -3.75%  "
$ws.Range("D45").Value = "1.255"
$ws.Range("E45").Value = "  -3.69%  "
$ws.Range("D46").Value = "2.271"
$ws.Range("E46").Value = "  +12.57%  "
$ws.Range("D47").Value = "0.5417"
$ws.Range("E47").Value = "  -3.50%  "
$ws.Range("D48").Value = "11.81"
$ws.Range("E48").Value = "  -3.00%  "
$ws.Range("D49").Value = "0.07146"
$ws.Range("E49").Value = "  +3.61%  "
$ws.Range("D50").Value = "1.879"
$ws.Range("E50").Value = "  -4.42%  "
$ws.Range("D51").Value = "110.05"
$ws.Range("E51").Value = "  -0.69%  "
